$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "67.273.35"
$ws.Range("E2").Value = "  +3.30%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.758.56"
$ws.Range("E3").Value = "  +7.19%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue "D5" "420.38"
$ws.Range("E5").Value = "  +0.46%  "

# Row 6 - Solana
Set-TextValue "D6" "131.84"
$ws.Range("E6").Value = "  -0.78%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.746.49"
$ws.Range("E7").Value = "  +7.22%  "

# Row 8 - XRP
Set-TextValue "D8" "0.650"
$ws.Range("E8").Value = "  -0.40%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.06%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.774"
$ws.Range("E10").Value = "  -0.79%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.186"
$ws.Range("E11").Value = "  +14.68%  "

# Row 12
Set-TextValue "D12" "0.0000422"
$ws.Range("E12").Value = "  +58.19%  "

# Row 13
Set-TextValue "D13" "42.91"
$ws.Range("E13").Value = "  -1.71%  "

# Row 14
Set-TextValue "D14" "10.40"
$ws.Range("E14").Value = "  +2.58%  "

# Row 15
Set-TextValue "D15" "4.351.97"
$ws.Range("E15").Value = "  +7.18%  "

# Row 16
$ws.Range("E16").Value = "  -0.83%  "

# Row 17
Set-TextValue "D17" "3.755.99"
$ws.Range("E17").Value = "  +6.96%  "

# Row 18
Set-TextValue "D18" "20.67"
$ws.Range("E18").Value = "  +0.34%  "

# Row 19
Set-TextValue "D19" "13.16"
$ws.Range("E19").Value = "  +2.74%  "

# Row 20
Set-TextValue "D20" "1.15"
$ws.Range("E20").Value = "  +2.72%  "

# Row 21
Set-TextValue "D21" "67.273.86"
$ws.Range("E21").Value = "  +3.51%  "

# Row 22
Set-TextValue "D22" "449.19"
$ws.Range("E22").Value = "  -1.03%  "

# Row 23
Set-TextValue "D23" "15.91"
$ws.Range("E23").Value = "  +19.12%  "

# Row 24
Set-TextValue "D24" "89.33"
$ws.Range("E24").Value = "  -1.14%  "

# Row 25
Set-TextValue "D25" "3.10"
$ws.Range("E25").Value = "  -3.63%  "

# Row 26
Set-TextValue "D26" "38.75"
$ws.Range("E26").Value = "  +13.35%  "

# Row 27
$ws.Range("E27").Value = "  -3.08%  "

# Row 28
$ws.Range("E28").Value = "  +1.63%  "

# Row 29
Set-TextValue "D29" "5.10"
$ws.Range("E29").Value = "  +5.55%  "

# Row 30
$ws.Range("E30").Value = "  +7.10%  "

# Row 31
Set-TextValue "D31" "12.70"
$ws.Range("E31").Value = "  +0.35%  "

# Row 32
Set-TextValue "D32" "2.70"
$ws.Range("E32").Value = "  -1.17%  "

# Row 33
Set-TextValue "D33" "7.29"
$ws.Range("E33").Value = "  -3.07%  "

# Row 34
$ws.Range("E34").Value = "  +1.51%  "

# Row 35
Set-TextValue "D35" "42.03"
$ws.Range("E35").Value = "  +5.08%  "

# Row 36
Set-TextValue "D36" "56.96"
$ws.Range("E36").Value = "  -0.27%  "

# Row 37
Set-TextValue "D37" "1.00"
$ws.Range("E37").Value = "  +0.11%  "

# Row 38
Set-TextValue "D38" "0.0493"
$ws.Range("E38").Value = "  -3.43%  "

# Row 39
Set-TextValue "D39" "0.0₃0783"
$ws.Range("E39").Value = "  +5.35%  "

# Row 40
$ws.Range("E40").Value = "  +0.18%  "

# Row 41 - ThetaToken
Set-TextValue "D41" "2.94"
$ws.Range("E41").Value = "  +25.82%  "

# Row 42 - now FirstDigitalUSD (was EnergySwap)
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D42" "0.997"
$ws.Range("E42").Value = "  -0.10%  "

# Row 43 - now EnergySwap (was FirstDigitalUSD)
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D43" "28.17"
$ws.Range("E43").Value = "  +30.04%  "

# Row 44 - LidoDAOToken
Set-TextValue "D44" "3.44"
$ws.Range("E44").Value = "  +3.59%  "

# Row 45 - now Monero (was ARBITRUM)
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D45" "148.82"
$ws.Range("E45").Value = "  +1.81%  "

# Row 46 - now ARBITRUM (was Monero)
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D46" "2.15"
$ws.Range("E46").Value = "  +6.94%  "

# Row 47 - ApeXProtocol
Set-TextValue "D47" "3.18"
$ws.Range("E47").Value = "  +23.87%  "

# Row 48
$ws.Range("E48").Value = "  -5.48%  "

# Row 49
Set-TextValue "D49" "2.65"
$ws.Range("E49").Value = "  -3.64%  "

# Row 50
Set-TextValue "D50" "4.34"
$ws.Range("E50").Value = "  -4.85%  "

# Row 51
Set-TextValue "D51" "0.308"
$ws.Range("E51").Value = "  -2.00%  "
